$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "2018 LEAVE CREDITS" (Table15) - add new monthly VL(5-0-0)
# accrual rows (1.25 earned per month) for Jun-Nov 2023.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("2018 LEAVE CREDITS")

$ws1.Range("C84").Value = 1.25

$ws1.Range("A85").Value = 45108
$ws1.Range("C85").Value = 1.25

$ws1.Range("A86").Value = 45139
$ws1.Range("C86").Value = 1.25

$ws1.Range("A87").Value = 45170
$ws1.Range("C87").Value = 1.25

$ws1.Range("A88").Value = 45200
$ws1.Range("C88").Value = 1.25

$ws1.Range("A89").Value = 45231
$ws1.Range("C89").Value = 1.25

# ---------------------------------------------------------------------
# Sheet "2017 LEAVE BALANCE" (Table1) - record leave usage for
# Oct-Dec 2023 (entered in the same order the author used, so new
# shared strings land on the same indexes as the source workbook).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2017 LEAVE BALANCE")

# Row 19 - SL(3-0-0), 10/23-25/2023
$ws2.Range("B19").Value = "SL(3-0-0)"
$ws2.Range("H19").Value = 3
$ws2.Range("K19").Value = "10/23-25/2023"

# Row 18 - VL(5-0-0), 11/14-17,20/2023
$ws2.Range("A18").Value = 45200
$ws2.Range("B18").Value = "VL(5-0-0)"
$ws2.Range("D18").Value = 5
$ws2.Range("K18").Value = "11/14-17,20/2023"

# Row 21 - SP(1-0-0), dated 12/15/2023 (45275)
$ws2.Range("B21").Value = "SP(1-0-0)"
$ws2.Range("K21").Value = 45275
$ws2.Range("K11").Copy()
$ws2.Range("K21").PasteSpecial(-4122)

# Row 20 - VL(4-0-0), 12/11-13, 18/2023
$ws2.Range("A20").Value = 45231
$ws2.Range("B20").Value = "VL(4-0-0)"
$ws2.Range("D20").Value = 4
$ws2.Range("K20").Value = "12/11-13, 18/2023"

# Row 22 - VL(4-0-0), 12/19,27-29/2023
$ws2.Range("B22").Value = "VL(4-0-0)"
$ws2.Range("D22").Value = 4
$ws2.Range("K22").Value = "12/19,27-29/2023"

# Row 23 - SL(1-0-0), dated 12/1/2023 (45261)
$ws2.Range("B23").Value = "SL(1-0-0)"
$ws2.Range("H23").Value = 1
$ws2.Range("K23").Value = 45261
$ws2.Range("K11").Copy()
$ws2.Range("K23").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Leave the UI pointed at the last-edited sheet/cell, matching the
# author's final on-screen selection.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("K22").Select()
